# Update the "想去人数" (interest count) figures in column F across the
# three sheets that carry event data, matching a newer scrape snapshot.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value  = 6
$wsExhibit.Range("F3").Value  = 12824
$wsExhibit.Range("F6").Value  = 60
$wsExhibit.Range("F10").Value = 12735
$wsExhibit.Range("F11").Value = 272
$wsExhibit.Range("F12").Value = 22
$wsExhibit.Range("F13").Value = 8648
$wsExhibit.Range("F14").Value = 7639
$wsExhibit.Range("F23").Value = 183
$wsExhibit.Range("F24").Value = 13

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F3").Value = 1

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value  = 6
$wsAll.Range("F4").Value  = 12824
$wsAll.Range("F7").Value  = 60
$wsAll.Range("F11").Value = 12735
$wsAll.Range("F12").Value = 272
$wsAll.Range("F13").Value = 22
$wsAll.Range("F14").Value = 8648
$wsAll.Range("F15").Value = 7639
$wsAll.Range("F23").Value = 1
$wsAll.Range("F25").Value = 183
$wsAll.Range("F26").Value = 13
